$wb = $excel.ActiveWorkbook

# --- Select the source range on Aggregate-type first so the post-copy
#     selection on that sheet ends up as a plain A1:N2 range-select
#     (matches the target diff, which updates Aggregate-type's saved
#     selection and drops its tabSelected flag in favour of the new sheet).
$src = $wb.Worksheets.Item("Aggregate-type")
$src.Range("A1:N2").Select()

# --- Duplicate "Aggregate-type" to the end of the workbook; this becomes
#     the new "PDF-Verifier" sheet (Excel's usual way of adding a sheet
#     that already carries the same layout/styles as a sibling sheet).
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "PDF-Verifier"

# --- Insert the new "Type" column as column B (shifts everything right).
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").Value = "Type"
$ws.Range("B2").Value = "REST"

# --- Fill in row 2 with the PDF-Verifier test case content. Values are
#     set in the same order the author's shared-string table grows in,
#     so newly interned strings land at the expected indices.
$ws.Range("C2").Value = "PDF-VERIFIER"
$ws.Range("D2").Value = "pdf"
$ws.Range("F2").Value = "application/pdf"
$ws.Range("O2").Value = " @pdf"
$ws.Range("J2").Value = "GET"
$ws.Range("A2").Value = "PDF-VERIFIER-2"
$ws.Range("H2").Value = "VirtualanStdType=PDF_TEMPLATE"
$ws.Range("E2").Value = "https://localhost:3000/api/pdf-validation.pdf"
$ws.Range("I2").Value = '{"from":["FROM:","TC01-01-API","APRIL CSR 2022","100000 FedExGndDrivefourthfloo","West Wing","Pittsburgh, PA 15108"],"to":["TO:","FEDEX ECONOMY LABEL VALIDATION","MICHAEL WESTEN","ourthfloo 528 NW 7TH AVENUE","APT. #2","MIAMI, FL 33136"]}'
$ws.Range("G2").Value = ""
$ws.Range("K2").Value = 200
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

$ws.Rows.Item(2).RowHeight = 166.5

# --- Drop the hyperlink that was carried over from the source sheet
#     (it pointed at Aggregate-type's Resource-column URL, which no
#     longer applies once this row's own data is in place).
$ws.Range("D2").Hyperlinks.Delete()

# --- Hyperlink on the URL cell, matching the pattern used on the other
#     sheets.
$pdfUrl = "https://localhost:3000/api/pdf-validation.pdf"
$ws.Hyperlinks.Add($ws.Range("E2"), $pdfUrl) | Out-Null
$ws.Range("E2").Value = $pdfUrl

# --- Re-create the data table over the new range (now A1:O2 with the
#     inserted Type column) using the same table style as the sibling
#     sheets' tables.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:O2"), $null, 1)
$lo.Name = "Table13"
$lo.TableStyle = "TableStyleMedium23"
$lo.ListColumns.Item("URL").DataBodyRange.Style = "Hyperlink"

# --- Final selection on the new sheet.
$ws.Range("H11:H25").Select()
